# Update "想去人数" (wish-to-go count) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1830
    $ws.Range("F5").Value = 1125
    $ws.Range("F6").Value = 1060
    $ws.Range("F9").Value = 97
}
